# Horarios actualizados Linea 141 - 47
# Scrape run at 01:22:42: each of the three schedule sheets gets its
# "Ultima actualizacion" / "Total filas" header refreshed and a new
# data row appended for the 215_ALUAR arrival observed at 01:22:42
# (02:58 elapsed => 96 minutes).

$wb = $excel.ActiveWorkbook

function Update-HorarioSheet {
    param(
        [string]$SheetName,
        [int]$NewRow,
        [int]$NewTotalFilas
    )

    $ws = $wb.Worksheets.Item($SheetName)

    $ws.Cells.Item(2, 1).Value = "Última actualización: 01:22:42"
    $ws.Cells.Item(3, 1).Value = "Total filas: $NewTotalFilas"

    $ws.Cells.Item($NewRow, 1).Value = "01:22:42"
    $ws.Cells.Item($NewRow, 2).Value = "02:58"
    $ws.Cells.Item($NewRow, 3).Value = "215_ALUAR"
    $ws.Cells.Item($NewRow, 4).Value = 96
}

# LP1912: rows 1-7 already used -> new row 8, total filas 2 -> 3
Update-HorarioSheet "LP1912" 8 3

# LP1912-215: rows 1-6 already used -> new row 7, total filas 1 -> 2
Update-HorarioSheet "LP1912-215" 7 2

# 6203-6173: rows 1-7 already used -> new row 8, total filas 2 -> 3
Update-HorarioSheet "6203-6173" 8 3
